$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Text"
$ws.Range("C1").Value = "Speaker"
$ws.Range("D1").Value = "Comments"

$ws.Range("A2").Value = "scene1_Scene1_Part1_S494"
$ws.Range("B2").Value = "This is a scene."
$ws.Range("C2").Value = "FRED"
$ws.Range("D2").Value = ""

$ws.Range("A3").Value = "scene1_Scene1_Part1_VXIU"
$ws.Range("B3").Value = "Part2"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""

$ws.Range("A4").Value = "scene1_Scene1_Part2_CF6W"
$ws.Range("B4").Value = "Part3"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""

$ws.Range("A5").Value = "scene1_Scene1_Part3_9MXL"
$ws.Range("B5").Value = "Let's see if this works, shall we?"
$ws.Range("C5").Value = "DAVE"
$ws.Range("D5").Value = ""

$ws.Range("A6").Value = "scene1_Scene1_Part3_YTMH"
$ws.Range("B6").Value = "Part4"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""

$ws.Range("A7").Value = "scene1_Scene1_Part4_T9GZ"
$ws.Range("B7").Value = "Go right"
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""

$ws.Range("A8").Value = "scene1_Scene1_Part4_F0PF"
$ws.Range("B8").Value = "Go right!"
$ws.Range("C8").Value = "FRED"
$ws.Range("D8").Value = ""

$ws.Range("A9").Value = "scene1_Scene1_Part4_9L7I"
$ws.Range("B9").Value = "Go left"
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""

$ws.Range("A10").Value = "scene1_Scene1_Part4_DNII"
$ws.Range("B10").Value = "Go left!"
$ws.Range("C10").Value = "FRED"
$ws.Range("D10").Value = ""

$ws.Range("A11").Value = "scene1_Scene1_Part4_Q8FK"
$ws.Range("B11").Value = "Skip it"
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""

$ws.Range("A12").Value = "scene1_Scene1_Part4_AJDP"
$ws.Range("B12").Value = "Carry on then."
$ws.Range("C12").Value = "FRED"
$ws.Range("D12").Value = ""

$ws.Range("A13").Value = "scene1_Scene1_Part4_0YY1"
$ws.Range("B13").Value = "Okay."
$ws.Range("C13").Value = "GEORGE"
$ws.Range("D13").Value = ""

$ws.Range("A14").Value = "scene1_Scene1_Part4_PZV1"
$ws.Range("B14").Value = "Back"
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""

$ws.Range("A15").Value = "scene1_Scene1_Right_WM69"
$ws.Range("B15").Value = "You sure you want to go right?"
$ws.Range("C15").Value = "GEORGE"
$ws.Range("D15").Value = ""

$ws.Range("A16").Value = "scene1_Scene1_Right_P8FP"
$ws.Range("B16").Value = "Back"
$ws.Range("C16").Value = ""
$ws.Range("D16").Value = ""

$ws.Range("A17").Value = "scene1_Scene1_Left_MIM6"
$ws.Range("B17").Value = "You sure you want to go left?"
$ws.Range("C17").Value = "GEORGE"
$ws.Range("D17").Value = ""

$ws.Range("A18").Value = "scene1_Scene1_Left_WXCN"
$ws.Range("B18").Value = "Back"
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""

$ws.Range("A19").Value = "scene1_Scene1_OtherContent_FSDK"
$ws.Range("B19").Value = "This content is nothing at all to do with Dink!"
$ws.Range("C19").Value = ""
$ws.Range("D19").Value = ""

$ws.Range("A20").Value = "scene1_OtherContent_VZWQ"
$ws.Range("B20").Value = "Back"
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = ""

$ws.Range("A21").Value = "main_Main_QU2R"
$ws.Range("B21").Value = "Barks"
$ws.Range("C21").Value = ""
$ws.Range("D21").Value = ""

$ws.Range("A22").Value = "main_Main_X20S"
$ws.Range("B22").Value = "Intro"
$ws.Range("C22").Value = ""
$ws.Range("D22").Value = ""

$ws.Range("A23").Value = "main_Main_OQ5O"
$ws.Range("B23").Value = "Intro2"
$ws.Range("C23").Value = ""
$ws.Range("D23").Value = ""

$ws.Range("A24").Value = "main_Main_NEAB"
$ws.Range("B24").Value = "TestScene"
$ws.Range("C24").Value = ""
$ws.Range("D24").Value = ""

$ws.Range("A25").Value = "main_Main_AD94"
$ws.Range("B25").Value = "Scene1"
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = ""

$ws.Range("A26").Value = "main_Intro_FDAP"
$ws.Range("B26").Value = "This is a test file."
$ws.Range("C26").Value = ""
$ws.Range("D26").Value = ""

$ws.Range("A27").Value = "main_Intro_EBU9"
$ws.Range("B27").Value = "Back"
$ws.Range("C27").Value = ""
$ws.Range("D27").Value = ""

$ws.Range("A28").Value = "main_Intro2_PCBU"
$ws.Range("B28").Value = "LAURA: This is an earlier line I am saying."
$ws.Range("C28").Value = ""
$ws.Range("D28").Value = ""

$ws.Range("A29").Value = "main_Intro2_QEUQ"
$ws.Range("B29").Value = "Back"
$ws.Range("C29").Value = ""
$ws.Range("D29").Value = ""

$ws.Range("A30").Value = "main_TestScene_16U4"
$ws.Range("B30").Value = "This is a line I am saying."
$ws.Range("C30").Value = "LAURA"
$ws.Range("D30").Value = ""

$ws.Range("A31").Value = "main_TestScene_FF1T"
$ws.Range("B31").Value = "This is another line."
$ws.Range("C31").Value = "LAURA"
$ws.Range("D31").Value = "VO: This comment goes to the voice actor., LOC: This comment goes to the localisers"

$ws.Range("A32").Value = "main_TestScene_BQ1E"
$ws.Range("B32").Value = "This is a loud line!"
$ws.Range("C32").Value = "FRED"
$ws.Range("D32").Value = ""

$ws.Range("A33").Value = "main_TestScene_IQIS"
$ws.Range("B33").Value = "Glad that's over with!"
$ws.Range("C33").Value = "FRED"
$ws.Range("D33").Value = ""

$ws.Range("A34").Value = "main_TestScene_MP0B"
$ws.Range("B34").Value = "Back"
$ws.Range("C34").Value = ""
$ws.Range("D34").Value = ""

$ws.Range("A35").Value = "main_Barks_O037"
$ws.Range("B35").Value = "Bark1"
$ws.Range("C35").Value = "FRED"
$ws.Range("D35").Value = ""

$ws.Range("A36").Value = "main_Barks_UWZ2"
$ws.Range("B36").Value = "Bark2"
$ws.Range("C36").Value = "FRED"
$ws.Range("D36").Value = ""

$ws.Range("A37").Value = "main_Barks_1ZG8"
$ws.Range("B37").Value = "Bark3"
$ws.Range("C37").Value = "FRED"
$ws.Range("D37").Value = ""

$ws.Range("A38").Value = "main_Barks_JFG1"
$ws.Range("B38").Value = "Bark4"
$ws.Range("C38").Value = "FRED"
$ws.Range("D38").Value = ""

$ws.Range("A39").Value = "main_Barks_X291"
$ws.Range("B39").Value = "Response to Bark 4."
$ws.Range("C39").Value = "JIM"
$ws.Range("D39").Value = ""

$ws.Range("A40").Value = "main_Barks_L2SX"
$ws.Range("B40").Value = "Bark5"
$ws.Range("C40").Value = "FRED"
$ws.Range("D40").Value = ""

$ws.Range("A41").Value = "main_Barks_N07F"
$ws.Range("B41").Value = "Bark6"
$ws.Range("C41").Value = "FRED"
$ws.Range("D41").Value = ""

$ws.Range("A42").Value = "main_Barks_83WH"
$ws.Range("B42").Value = "Back"
$ws.Range("C42").Value = ""
$ws.Range("D42").Value = ""

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D42"))
